$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix labeling on the Shadow War row: "Shadow War 7PM" -> "Shadow War 6PM"
$ws.Range("A12").Value = "Shadow War 6PM"

# Leave selection on the edited area, matching the saved cursor position
$ws.Range("A13").Select()
